$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 30455.889
$ws.Range("I10").Value = 104
$ws.Range("J10").Value = 34249.875
$ws.Range("K10").Value = 104
$ws.Range("L10").Value = 34249.875
$ws.Range("M10").Value = 189
$ws.Range("N10").Value = -34835.875
$ws.Range("H86").Value = 4273.778
$ws.Range("J86").Value = 5903.6
$ws.Range("L86").Value = 5903.6
$ws.Range("N86").Value = -8149.6
$ws.Range("H89").Value = 4273.778
$ws.Range("J89").Value = 5903.6
$ws.Range("L89").Value = 29518
$ws.Range("N89").Value = -40750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3972.5454
$ws.Range("I32").Value = 3494.7896
$ws.Range("K32").Value = 3494.7896
$ws.Range("M32").Value = -3207.7896
$ws.Range("H37").Value = 7857.143
$ws.Range("J37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("N37").Value = -10546
$ws.Range("H61").Value = 924
$ws.Range("I61").Value = 923.5
$ws.Range("K61").Value = 923.5
$ws.Range("M61").Value = -711.5
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H74").Value = 1155.3334
$ws.Range("I74").Value = 890.7778
$ws.Range("J74").Value = 1949
$ws.Range("K74").Value = 890.7778
$ws.Range("L74").Value = 1949
$ws.Range("M74").Value = -16.77779999999996
$ws.Range("N74").Value = -3697
$ws.Range("H77").Value = 1155.3334
$ws.Range("I77").Value = 890.7778
$ws.Range("J77").Value = 1949
$ws.Range("K77").Value = 4453.889
$ws.Range("L77").Value = 9745
$ws.Range("M77").Value = -85.88900000000012
$ws.Range("N77").Value = -18481
$ws.Range("H136").Value = 924
$ws.Range("I136").Value = 923.5
$ws.Range("K136").Value = 2770.5
$ws.Range("M136").Value = -220.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 20078.8
$ws.Range("I11").Value = 131.33333
$ws.Range("J11").Value = 50000
$ws.Range("K11").Value = 131.33333
$ws.Range("L11").Value = 50000
$ws.Range("M11").Value = 8.666670000000011
$ws.Range("N11").Value = -50280
$ws.Range("H94").Value = 1653.8235
$ws.Range("I94").Value = 1115.3572
$ws.Range("K94").Value = 1115.3572
$ws.Range("M94").Value = -664.3571999999999
$ws.Range("H134").Value = 5816.091
$ws.Range("I134").Value = 5447.7
$ws.Range("K134").Value = 16343.1
$ws.Range("M134").Value = -13808.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 481.42856
$ws.Range("H93").Value = 4407
$ws.Range("I93").Value = 4407
$ws.Range("K93").Value = 4407
$ws.Range("M93").Value = -2535
$ws.Range("H132").Value = 224.75
$ws.Range("I132").Value = 224.75
$ws.Range("K132").Value = 674.25
$ws.Range("M132").Value = 1855.75
$ws.Range("H134").Value = 2328.5
$ws.Range("I134").Value = 2650
$ws.Range("K134").Value = 7950
$ws.Range("M134").Value = -5415

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 71428744
$ws.Range("I4").Value = 189.53847
$ws.Range("K4").Value = 568.61541
$ws.Range("M4").Value = -456.61541
$ws.Range("H5").Value = 1161
$ws.Range("I5").Value = 1266.75
$ws.Range("J5").Value = 1100.5714
$ws.Range("K5").Value = 3800.25
$ws.Range("L5").Value = 3301.7142
$ws.Range("M5").Value = -3688.25
$ws.Range("N5").Value = -3525.7142
$ws.Range("H68").Value = 4397.1816
$ws.Range("I68").Value = 3278.4
$ws.Range("K68").Value = 9835.200000000001
$ws.Range("M68").Value = -9024.200000000001
$ws.Range("H71").Value = 4397.1816
$ws.Range("I71").Value = 3278.4
$ws.Range("K71").Value = 29505.6
$ws.Range("M71").Value = -25449.6
$ws.Range("H129").Value = 1202.1
$ws.Range("I129").Value = 752.625
$ws.Range("K129").Value = 2257.875
$ws.Range("M129").Value = 2742.125
$ws.Range("H131").Value = 919
$ws.Range("I131").Value = 647
$ws.Range("J131").Value = 957.8570999999999
$ws.Range("K131").Value = 1941
$ws.Range("L131").Value = 2873.5713
$ws.Range("M131").Value = 3099
$ws.Range("N131").Value = -12953.5713
$ws.Range("H135").Value = 1161
$ws.Range("I135").Value = 1266.75
$ws.Range("J135").Value = 1100.5714
$ws.Range("K135").Value = 11400.75
$ws.Range("L135").Value = 9905.142600000001
$ws.Range("M135").Value = -8865.75
$ws.Range("N135").Value = -14975.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 852.7143
$ws.Range("I102").Value = 852.7143
$ws.Range("K102").Value = 852.7143
$ws.Range("M102").Value = 769.2857
$ws.Range("H113").Value = 4784.2
$ws.Range("J113").Value = 4999
$ws.Range("L113").Value = 4999
$ws.Range("N113").Value = -9339
$ws.Range("H132").Value = 4005.7576
$ws.Range("I132").Value = 3840.577
$ws.Range("J132").Value = 4619.2856
$ws.Range("K132").Value = 11521.731
$ws.Range("L132").Value = 13857.8568
$ws.Range("M132").Value = -8991.731
$ws.Range("N132").Value = -18917.8568
$ws.Range("H138").Value = 120000
$ws.Range("J138").Value = 120000
$ws.Range("L138").Value = 120000
$ws.Range("N138").Value = -130280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3867.5881
$ws.Range("I46").Value = 3399.9
$ws.Range("J46").Value = 4535.7144
$ws.Range("K46").Value = 3399.9
$ws.Range("L46").Value = 4535.7144
$ws.Range("M46").Value = -3211.9
$ws.Range("N46").Value = -4911.7144
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = $null
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H61").Value = 1100.75
$ws.Range("I61").Value = 1100.75
$ws.Range("K61").Value = 1100.75
$ws.Range("M61").Value = -898.75
$ws.Range("H113").Value = 1100.75
$ws.Range("I113").Value = 1100.75
$ws.Range("K113").Value = 1100.75
$ws.Range("M113").Value = 1069.25
$ws.Range("H132").Value = 6434
$ws.Range("I132").Value = 3892.5
$ws.Range("K132").Value = 11677.5
$ws.Range("M132").Value = -9147.5
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5001751
$ws.Range("J3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("N3").Value = -2728
$ws.Range("H132").Value = 3933
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 16500
$ws.Range("M132").Value = -13970
$ws.Range("H136").Value = 1887.75
$ws.Range("I136").Value = 1887.75
$ws.Range("K136").Value = 5663.25
$ws.Range("M136").Value = -3113.25
